$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 180.27272
$ws.Range("I4").Value = 98.09999999999999
$ws.Range("K4").Value = 98.09999999999999
$ws.Range("M4").Value = 15.90000000000001
$ws.Range("H39").Value = 816.7692
$ws.Range("I39").Value = 56.18182
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 168.54546
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = 127.45454
$ws.Range("N39").Value = -15592
$ws.Range("H64").Value = 5907.393
$ws.Range("I64").Value = 5477.273
$ws.Range("J64").Value = 7484.5
$ws.Range("K64").Value = 5477.273
$ws.Range("L64").Value = 7484.5
$ws.Range("M64").Value = -5229.273
$ws.Range("N64").Value = -7980.5
$ws.Range("H67").Value = 5907.393
$ws.Range("I67").Value = 5477.273
$ws.Range("J67").Value = 7484.5
$ws.Range("K67").Value = 5477.273
$ws.Range("L67").Value = 7484.5
$ws.Range("M67").Value = -4619.273
$ws.Range("N67").Value = -9200.5
$ws.Range("H88").Value = 1978.3334
$ws.Range("J88").Value = 2100.25
$ws.Range("L88").Value = 2100.25
$ws.Range("N88").Value = -2912.25
$ws.Range("H91").Value = 1978.3334
$ws.Range("J91").Value = 2100.25
$ws.Range("L91").Value = 2100.25
$ws.Range("N91").Value = -4908.25
$ws.Range("H132").Value = 5039.4165
$ws.Range("I132").Value = 5226.1304
$ws.Range("J132").Value = 745
$ws.Range("K132").Value = 15678.3912
$ws.Range("L132").Value = 2235
$ws.Range("M132").Value = -13148.3912
$ws.Range("N132").Value = -7295
$ws.Range("H138").Value = 4583.5347
$ws.Range("I138").Value = 2092.3572
$ws.Range("J138").Value = 5786.1724
$ws.Range("K138").Value = 6277.071599999999
$ws.Range("L138").Value = 17358.5172
$ws.Range("M138").Value = -1137.071599999999
$ws.Range("N138").Value = -27638.5172

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2236.5789
$ws.Range("I2").Value = 1091.8182
$ws.Range("J2").Value = 3810.625
$ws.Range("K2").Value = 1091.8182
$ws.Range("L2").Value = 3810.625
$ws.Range("M2").Value = -978.8181999999999
$ws.Range("N2").Value = -4036.625
$ws.Range("H32").Value = 5575.3193
$ws.Range("I32").Value = 4937.0513
$ws.Range("K32").Value = 4937.0513
$ws.Range("M32").Value = -4650.0513
$ws.Range("H63").Value = 2534.8333
$ws.Range("I63").Value = 2441.8
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2441.8
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1755.8
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2534.8333
$ws.Range("I66").Value = 2441.8
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 12209
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -8777
$ws.Range("N66").Value = -21864
$ws.Range("H116").Value = 2236.5789
$ws.Range("I116").Value = 1091.8182
$ws.Range("J116").Value = 3810.625
$ws.Range("K116").Value = 1091.8182
$ws.Range("L116").Value = 3810.625
$ws.Range("M116").Value = 1202.1818
$ws.Range("N116").Value = -8398.625
$ws.Range("H132").Value = 28615408
$ws.Range("I132").Value = 4326.8965
$ws.Range("K132").Value = 12980.6895
$ws.Range("M132").Value = -10450.6895

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2236.5789
$ws.Range("I3").Value = 1091.8182
$ws.Range("J3").Value = 3810.625
$ws.Range("K3").Value = 1091.8182
$ws.Range("L3").Value = 3810.625
$ws.Range("M3").Value = -977.8181999999999
$ws.Range("N3").Value = -4038.625
$ws.Range("H20").Value = 2205.1562
$ws.Range("I20").Value = 1560.381
$ws.Range("J20").Value = 3436.0908
$ws.Range("K20").Value = 1560.381
$ws.Range("L20").Value = 3436.0908
$ws.Range("M20").Value = -1313.381
$ws.Range("N20").Value = -3930.0908
$ws.Range("H44").Value = 34999
$ws.Range("I44").Value = 34999
$ws.Range("K44").Value = 34999
$ws.Range("M44").Value = -34502
$ws.Range("H96").Value = 39303.668
$ws.Range("I96").Value = 8127.1665
$ws.Range("J96").Value = 101656.664
$ws.Range("K96").Value = 8127.1665
$ws.Range("L96").Value = 101656.664
$ws.Range("M96").Value = -5381.1665
$ws.Range("N96").Value = -107148.664

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("M40").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3666.25
$ws.Range("I5").Value = 1776.6666
$ws.Range("J5").Value = 4800
$ws.Range("K5").Value = 5329.9998
$ws.Range("L5").Value = 14400
$ws.Range("M5").Value = -5217.9998
$ws.Range("N5").Value = -14624
$ws.Range("H40").Value = 174
$ws.Range("I40").Value = 62.5
$ws.Range("J40").Value = 307.8
$ws.Range("K40").Value = 250
$ws.Range("L40").Value = 1231.2
$ws.Range("M40").Value = -181
$ws.Range("N40").Value = -1369.2
$ws.Range("H56").Value = 21888.309
$ws.Range("I56").Value = 21888.309
$ws.Range("K56").Value = 21888.309
$ws.Range("M56").Value = -21358.309
$ws.Range("H113").Value = 2575.6667
$ws.Range("I113").Value = 233
$ws.Range("J113").Value = 3747
$ws.Range("K113").Value = 699
$ws.Range("L113").Value = 11241
$ws.Range("M113").Value = 1471
$ws.Range("N113").Value = -15581
$ws.Range("H122").Value = 2477.4614
$ws.Range("J122").Value = 3344.6667
$ws.Range("L122").Value = 30102.0003
$ws.Range("N122").Value = -35002.0003
$ws.Range("H135").Value = 3666.25
$ws.Range("I135").Value = 1776.6666
$ws.Range("J135").Value = 4800
$ws.Range("K135").Value = 15989.9994
$ws.Range("L135").Value = 43200
$ws.Range("M135").Value = -13454.9994
$ws.Range("N135").Value = -48270

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1048.32
$ws.Range("I97").Value = 1141.8334
$ws.Range("J97").Value = 807.8570999999999
$ws.Range("K97").Value = 1141.8334
$ws.Range("L97").Value = 807.8570999999999
$ws.Range("M97").Value = -645.8334
$ws.Range("N97").Value = -1799.8571
$ws.Range("H122").Value = 38463156
$ws.Range("I122").Value = 1531.6666
$ws.Range("J122").Value = 125001810
$ws.Range("K122").Value = 4594.9998
$ws.Range("L122").Value = 375005430
$ws.Range("M122").Value = -2144.9998
$ws.Range("N122").Value = -375010330
$ws.Range("H132").Value = 2109.7144
$ws.Range("I132").Value = 2112.4783
$ws.Range("K132").Value = 6337.4349
$ws.Range("M132").Value = -3807.4349

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 859.4783
$ws.Range("I55").Value = 517.25
$ws.Range("J55").Value = 1232.8182
$ws.Range("K55").Value = 517.25
$ws.Range("L55").Value = 1232.8182
$ws.Range("M55").Value = -344.25
$ws.Range("N55").Value = -1578.8182
$ws.Range("H93").Value = 660967.8
$ws.Range("I93").Value = 1526.3529
$ws.Range("J93").Value = 2529385.2
$ws.Range("K93").Value = 1526.3529
$ws.Range("L93").Value = 2529385.2
$ws.Range("M93").Value = -278.3529000000001
$ws.Range("N93").Value = -2531881.2
$ws.Range("H132").Value = 6942.3184
$ws.Range("I132").Value = 2913.2104
$ws.Range("J132").Value = 32460
$ws.Range("K132").Value = 8739.6312
$ws.Range("L132").Value = 97380
$ws.Range("M132").Value = -6209.6312
$ws.Range("N132").Value = -102440
$ws.Range("H136").Value = 5220.25
$ws.Range("J136").Value = 11184.2
$ws.Range("L136").Value = 33552.60000000001
$ws.Range("N136").Value = -38652.60000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2425.4546
$ws.Range("I81").Value = 1131
$ws.Range("J81").Value = 3978.8
$ws.Range("K81").Value = 2262
$ws.Range("L81").Value = 7957.6
$ws.Range("M81").Value = -1201
$ws.Range("N81").Value = -10079.6
$ws.Range("H84").Value = 2425.4546
$ws.Range("I84").Value = 1131
$ws.Range("J84").Value = 3978.8
$ws.Range("K84").Value = 11310
$ws.Range("L84").Value = 39788
$ws.Range("M84").Value = -6006
$ws.Range("N84").Value = -50396
$ws.Range("H100").Value = 55556664
$ws.Range("I100").Value = 58824468
$ws.Range("K100").Value = 117648936
$ws.Range("M100").Value = -117648395
$ws.Range("H113").Value = 1721.375
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 6000
$ws.Range("N113").Value = -10340
